$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "maj template comment à la fin" - move the Comment column to the end
# (column K), putting Result in its former place (column J). This swaps
# the content of columns J and K across the header + the three
# description rows below it.

# Row 1 (headers)
$ws.Range("J1").Value = "Result"
$ws.Range("K1").Value = "Comment"

# Row 2 (field descriptions)
$ws.Range("K2").Value = "# Commentaire"
$ws.Range("J2").Value = ""

# Row 3 (field types)
$ws.Range("J3").Value = "#float"
$ws.Range("K3").Value = "#string"

# Row 4 (format notes)
$ws.Range("K4").Value = "# format: texte libre"
$ws.Range("J4").Value = ""

# Row 5 (examples) - both columns already blank, nothing to change.
